$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.204545672665859
$ws.Range("C2").Value = 0.3310472128250979
$ws.Range("E2").Value = 0.8703704342090646
$ws.Range("F2").Value = 1.903892237169117
$ws.Range("G2").Value = 0.1744859203586913
$ws.Range("H2").Value = 0.3633809362774727
$ws.Range("J2").Value = 0.01909735548007241
$ws.Range("O2").Value = 0.9797718923104668

$ws.Range("B3").Value = 1.054809573213731
$ws.Range("C3").Value = 0.292662111559963
$ws.Range("E3").Value = 0.8380285508793719
$ws.Range("F3").Value = 1.877963039739072
$ws.Range("G3").Value = 0.1793273158311735
$ws.Range("H3").Value = 0.3702880574768272
$ws.Range("J3").Value = 0.01958290440316546
$ws.Range("O3").Value = 1.004227446932546

$ws.Range("B4").Value = 0.9625168390858789
$ws.Range("C4").Value = 0.2690063626012602
$ws.Range("E4").Value = 0.8185054343311435
$ws.Range("F4").Value = 1.863390112672079
$ws.Range("G4").Value = 0.1826288461569199
$ws.Range("H4").Value = 0.3748290706111739
$ws.Range("J4").Value = 0.01990684410387011
$ws.Range("O4").Value = 1.020559826869551

$ws.Range("B5").Value = 0.9248203684108489
$ws.Range("C5").Value = 0.259345322087114
$ws.Range("E5").Value = 0.8106339367512874
$ws.Range("F5").Value = 1.857789503322863
$ws.Range("G5").Value = 0.1840564276807868
$ws.Range("H5").Value = 0.3767548545397936
$ws.Range("J5").Value = 0.02004531226641681
$ws.Range("O5").Value = 1.027545027140569

$ws.Range("B6").Value = 0.918555755730722
$ws.Range("C6").Value = 0.2577398606007364
$ws.Range("E6").Value = 0.8093319773342387
$ws.Range("F6").Value = 1.856879912728857
$ws.Range("G6").Value = 0.1842984247217814
$ws.Range("H6").Value = 0.3770791718862938
$ws.Range("J6").Value = 0.02006869408859568
$ws.Range("O6").Value = 1.028724780487188

$ws.Range("B7").Value = 0.962008797758358
$ws.Range("C7").Value = 0.2688761549438539
$ws.Range("E7").Value = 0.818398934961067
$ws.Range("F7").Value = 1.863313213656582
$ws.Range("G7").Value = 0.1826477669655944
$ws.Range("H7").Value = 0.3748547377925426
$ws.Range("J7").Value = 0.01990868541694812
$ws.Range("O7").Value = 1.020652699006703

$ws.Range("B8").Value = 1.152991693087301
$ws.Range("C8").Value = 0.3178305707329798
$ws.Range("E8").Value = 0.8591494909342003
$ws.Range("F8").Value = 1.894671568879204
$ws.Range("G8").Value = 0.1760866586009584
$ws.Range("H8").Value = 0.365700150698693
$ws.Range("J8").Value = 0.01925939763055595
$ws.Range("O8").Value = 0.9879300511477709

$ws.Range("B9").Value = 1.524601862754253
$ws.Range("C9").Value = 0.4131091233880966
$ws.Range("E9").Value = 0.9417190015913803
$ws.Range("F9").Value = 1.966905970812576
$ws.Range("G9").Value = 0.1658532310371683
$ws.Range("H9").Value = 0.3501353052020661
$ws.Range("J9").Value = 0.0181922836278634
$ws.Range("O9").Value = 0.9342712453734521

$ws.Range("B10").Value = 1.795749640745953
$ws.Range("C10").Value = 0.4826383230123383
$ws.Range("E10").Value = 1.004011011124248
$ws.Range("F10").Value = 2.0265987264508
$ws.Range("G10").Value = 0.1599715912918427
$ws.Range("H10").Value = 0.3401642158149372
$ws.Range("J10").Value = 0.0175357671374492
$ws.Range("O10").Value = 0.901342121638919

$ws.Range("B11").Value = 1.918675295204991
$ws.Range("C11").Value = 0.5141601698314844
$ws.Range("E11").Value = 1.032705054447689
$ws.Range("F11").Value = 2.055209040224241
$ws.Range("G11").Value = 0.1576580973278254
$ws.Range("H11").Value = 0.3359480237190837
$ws.Range("J11").Value = 0.01726517554624252
$ws.Range("O11").Value = 0.8877904354331463

$ws.Range("B12").Value = 1.965161345256433
$ws.Range("C12").Value = 0.5260805737848955
$ws.Range("E12").Value = 1.043622123624502
$ws.Range("F12").Value = 2.066253534911112
$ws.Range("G12").Value = 0.1568346623856698
$ws.Range("H12").Value = 0.3343976037244687
$ws.Range("J12").Value = 0.01716677729585747
$ws.Range("O12").Value = 0.8828656389110705

$ws.Range("B13").Value = 1.955152597536994
$ws.Range("C13").Value = 0.5235140362239576
$ws.Range("E13").Value = 1.041268659155151
$ws.Range("F13").Value = 2.063865531441309
$ws.Range("G13").Value = 0.1570096538725991
$ws.Range("H13").Value = 0.334729458434559
$ws.Range("J13").Value = 0.01718778761571471
$ws.Range("O13").Value = 0.883917053157802

$ws.Range("B14").Value = 1.922501013113731
$ws.Range("C14").Value = 0.5151411985685854
$ws.Range("E14").Value = 1.033602181527129
$ws.Range("F14").Value = 2.056113454889044
$ws.Range("G14").Value = 0.1575892947159332
$ws.Range("H14").Value = 0.3358195435377098
$ws.Range("J14").Value = 0.01725699852222107
$ws.Range("O14").Value = 0.8873811130770974

$ws.Range("B15").Value = 1.902492660828216
$ws.Range("C15").Value = 0.5100104558345606
$ws.Range("E15").Value = 1.028912915238493
$ws.Range("F15").Value = 2.051392514390187
$ws.Range("G15").Value = 0.1579512128889249
$ws.Range("H15").Value = 0.3364932691423235
$ws.Range("J15").Value = 0.01729992309607198
$ws.Range("O15").Value = 0.8895299459295103

$ws.Range("B16").Value = 1.787707459327748
$ws.Range("C16").Value = 0.4805760676629234
$ws.Range("E16").Value = 1.002142965582934
$ws.Range("F16").Value = 2.024758358844792
$ws.Range("G16").Value = 0.1601301182594312
$ws.Range("H16").Value = 0.3404462031860049
$ws.Range("J16").Value = 0.01755401848070193
$ws.Range("O16").Value = 0.9022566265633287

$ws.Range("B17").Value = 1.717180730984182
$ws.Range("C17").Value = 0.4624909599609737
$ws.Range("E17").Value = 0.9858118620043825
$ws.Range("F17").Value = 2.008792772987022
$ws.Range("G17").Value = 0.161559950406172
$ws.Range("H17").Value = 0.3429532177930525
$ws.Range("J17").Value = 0.01771710890150757
$ws.Range("O17").Value = 0.9104308994113666

$ws.Range("B18").Value = 1.676576171772581
$ws.Range("C18").Value = 0.4520788287088635
$ws.Range("E18").Value = 0.976452259885491
$ws.Range("F18").Value = 1.999746752029992
$ws.Range("G18").Value = 0.1624164002156832
$ws.Range("H18").Value = 0.3444252673904842
$ws.Range("J18").Value = 0.01781355341823776
$ws.Range("O18").Value = 0.9152668338710299

$ws.Range("B19").Value = 1.662821473250801
$ws.Range("C19").Value = 0.4485517574687492
$ws.Range("E19").Value = 0.9732890360311188
$ws.Range("F19").Value = 1.996707416167396
$ws.Range("G19").Value = 0.1627122099321241
$ws.Range("H19").Value = 0.3449288388516578
$ws.Range("J19").Value = 0.01784666013428726
$ws.Range("O19").Value = 0.9169272171594969

$ws.Range("B20").Value = 1.724692524340639
$ws.Range("C20").Value = 0.464417197104467
$ws.Range("E20").Value = 0.9875468573412718
$ws.Range("F20").Value = 2.010478154231777
$ws.Range("G20").Value = 0.161404214367181
$ws.Range("H20").Value = 0.3426832271701628
$ws.Range("J20").Value = 0.01769947425378149
$ws.Range("O20").Value = 0.909546822058573

$ws.Range("B21").Value = 1.932093312376935
$ws.Range("C21").Value = 0.5176009509015671
$ws.Range("E21").Value = 1.035852621432895
$ws.Range("F21").Value = 2.058384709368966
$ws.Range("G21").Value = 0.1574176071430387
$ws.Range("H21").Value = 0.335498104847133
$ws.Range("J21").Value = 0.01723655890319087
$ws.Range("O21").Value = 0.8863580058935838

$ws.Range("B22").Value = 2.067271942997536
$ws.Range("C22").Value = 0.5522646120670629
$ws.Range("E22").Value = 1.067721961016034
$ws.Range("F22").Value = 2.090921216561355
$ws.Range("G22").Value = 0.1551191878904987
$ws.Range("H22").Value = 0.3310713572468771
$ws.Range("J22").Value = 0.01695774937481076
$ws.Range("O22").Value = 0.8724096766521967

$ws.Range("B23").Value = 1.995159285939906
$ws.Range("C23").Value = 0.5337729178770587
$ws.Range("E23").Value = 1.050685389903762
$ws.Range("F23").Value = 2.073443269219297
$ws.Range("G23").Value = 0.1563176147853866
$ws.Range("H23").Value = 0.3334093071608066
$ws.Range("J23").Value = 0.01710437258992492
$ws.Range("O23").Value = 0.8797432038360284

$ws.Range("B24").Value = 1.72129662241241
$ws.Range("C24").Value = 0.4635463912151181
$ws.Range("E24").Value = 0.9867623744684977
$ws.Range("F24").Value = 2.009715779758238
$ws.Range("G24").Value = 0.1614745154559714
$ws.Range("H24").Value = 0.3428051942092978
$ws.Range("J24").Value = 0.01770743852746648
$ws.Range("O24").Value = 0.9099460885550315

$ws.Range("B25").Value = 1.424393342791575
$ws.Range("C25").Value = 0.3874143241890238
$ws.Range("E25").Value = 0.9190964260263428
$ws.Range("F25").Value = 1.946207118618588
$ws.Range("G25").Value = 0.1683364858298262
$ws.Range("H25").Value = 0.3540894981431535
$ws.Range("J25").Value = 0.01845870917227899
$ws.Range("O25").Value = 0.9476532121754957
